$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 454.42856
$ws.Range("I107").Value = 505.16666
$ws.Range("J107").Value = 150
$ws.Range("K107").Value = 505.16666
$ws.Range("L107").Value = 150
$ws.Range("M107").Value = 1414.83334
$ws.Range("N107").Value = -3990

$ws.Range("H112").Value = 6117.6665
$ws.Range("J112").Value = 3897.5789
$ws.Range("L112").Value = 11692.7367
$ws.Range("N112").Value = -13908.7367

$ws.Range("H113").Value = 18193462
$ws.Range("J113").Value = 5999
$ws.Range("L113").Value = 5999
$ws.Range("N113").Value = -12507

$ws.Range("H132").Value = 6814.019
$ws.Range("I132").Value = 5501.0684
$ws.Range("K132").Value = 16503.2052
$ws.Range("M132").Value = -13973.2052

$ws.Range("H136").Value = 69773.2
$ws.Range("J136").Value = 69773.2
$ws.Range("L136").Value = 69773.2
$ws.Range("N136").Value = -79973.2

$ws.Range("H137").Value = 12736.131
$ws.Range("I137").Value = 6158.6
$ws.Range("J137").Value = 17795.77
$ws.Range("K137").Value = 18475.8
$ws.Range("L137").Value = 53387.31
$ws.Range("M137").Value = -15925.8
$ws.Range("N137").Value = -58487.31

$ws.Range("H141").Value = 774.1667
$ws.Range("I141").Value = 784.75
$ws.Range("K141").Value = 2354.25
$ws.Range("M141").Value = 2825.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 11557.2
$ws.Range("J2").Value = 25944.5
$ws.Range("L2").Value = 25944.5
$ws.Range("N2").Value = -26170.5

$ws.Range("H5").Value = 175.11111
$ws.Range("I5").Value = 155.2
$ws.Range("K5").Value = 155.2
$ws.Range("M5").Value = -43.19999999999999

$ws.Range("H25").Value = 5503.6
$ws.Range("I25").Value = 1875
$ws.Range("J25").Value = 20018
$ws.Range("K25").Value = 1875
$ws.Range("L25").Value = 20018
$ws.Range("M25").Value = -1473
$ws.Range("N25").Value = -20822

$ws.Range("H32").Value = 7553.58
$ws.Range("I32").Value = 1374.8
$ws.Range("J32").Value = 21970.732
$ws.Range("K32").Value = 1374.8
$ws.Range("L32").Value = 21970.732
$ws.Range("M32").Value = -1087.8
$ws.Range("N32").Value = -22544.732

$ws.Range("H74").Value = 11054.75
$ws.Range("I74").Value = 2633.6667
$ws.Range("K74").Value = 2633.6667
$ws.Range("M74").Value = -1759.6667

$ws.Range("H77").Value = 11054.75
$ws.Range("I77").Value = 2633.6667
$ws.Range("K77").Value = 13168.3335
$ws.Range("M77").Value = -8800.333500000001

$ws.Range("H116").Value = 11557.2
$ws.Range("J116").Value = 25944.5
$ws.Range("L116").Value = 25944.5
$ws.Range("N116").Value = -30532.5

$ws.Range("H135").Value = 128734.5
$ws.Range("J135").Value = 128734.5
$ws.Range("L135").Value = 128734.5
$ws.Range("N135").Value = -138874.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 11557.2
$ws.Range("J3").Value = 25944.5
$ws.Range("L3").Value = 25944.5
$ws.Range("N3").Value = -26172.5

$ws.Range("H4").Value = 175.11111
$ws.Range("I4").Value = 155.2
$ws.Range("K4").Value = 155.2
$ws.Range("M4").Value = -40.19999999999999

$ws.Range("H25").Value = 3271.1428
$ws.Range("I25").Value = 1779.6
$ws.Range("K25").Value = 1779.6
$ws.Range("M25").Value = -1544.6

$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

$ws.Range("H105").Value = 1326.7931
$ws.Range("I105").Value = 661
$ws.Range("K105").Value = 661
$ws.Range("M105").Value = 1086

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14338.556
$ws.Range("I31").Value = 11558.454
$ws.Range("K31").Value = 11558.454
$ws.Range("M31").Value = -11263.454

$ws.Range("H34").Value = 14338.556
$ws.Range("I34").Value = 11558.454
$ws.Range("K34").Value = 11558.454
$ws.Range("M34").Value = -11356.454

$ws.Range("H86").Value = 9474.941000000001
$ws.Range("I86").Value = 11437.777
$ws.Range("J86").Value = 7266.75
$ws.Range("K86").Value = 11437.777
$ws.Range("L86").Value = 7266.75
$ws.Range("M86").Value = -10314.777
$ws.Range("N86").Value = -9512.75

$ws.Range("H89").Value = 9474.941000000001
$ws.Range("I89").Value = 11437.777
$ws.Range("J89").Value = 7266.75
$ws.Range("K89").Value = 57188.885
$ws.Range("L89").Value = 36333.75
$ws.Range("M89").Value = -51572.885
$ws.Range("N89").Value = -47565.75

$ws.Range("H110").Value = 79995
$ws.Range("J110").Value = 79995
$ws.Range("L110").Value = 79995
$ws.Range("N110").Value = -88175

$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").ClearContents()
$ws.Range("N111").Value = 0

$ws.Range("H116").Value = 75750
$ws.Range("J116").Value = 75750
$ws.Range("L116").Value = 75750

$ws.Range("H132").Value = 6260.607
$ws.Range("I132").Value = 1980.625
$ws.Range("K132").Value = 5941.875
$ws.Range("M132").Value = -3411.875

$ws.Range("H134").Value = 25647060
$ws.Range("I134").Value = 1938.579
$ws.Range("J134").Value = 50009924
$ws.Range("K134").Value = 5815.737
$ws.Range("L134").Value = 150029772
$ws.Range("M134").Value = -3280.737
$ws.Range("N134").Value = -150034842

$ws.Range("H135").Value = 73082.71000000001
$ws.Range("J135").Value = 73082.71000000001
$ws.Range("L135").Value = 73082.71000000001
$ws.Range("N135").Value = -83222.71000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 253.375
$ws.Range("I14").Value = 253.375
$ws.Range("K14").Value = 760.125
$ws.Range("M14").Value = -587.125

$ws.Range("H23").Value = 505
$ws.Range("I23").Value = 360.625
$ws.Range("J23").Value = 593.8461
$ws.Range("K23").Value = 1081.875
$ws.Range("L23").Value = 1781.5383
$ws.Range("M23").Value = -846.875
$ws.Range("N23").Value = -2251.5383

$ws.Range("H81").Value = 3574649.5
$ws.Range("I81").Value = 1523.5
$ws.Range("K81").Value = 4570.5
$ws.Range("M81").Value = -3447.5

$ws.Range("H84").Value = 3574649.5
$ws.Range("I84").Value = 1523.5
$ws.Range("K84").Value = 13711.5
$ws.Range("M84").Value = -8095.5

$ws.Range("H103").Value = 4943.9375
$ws.Range("I103").Value = 311
$ws.Range("J103").Value = 7049.8184
$ws.Range("K103").Value = 933
$ws.Range("L103").Value = 21149.4552
$ws.Range("M103").Value = -54
$ws.Range("N103").Value = -22907.4552

$ws.Range("H119").Value = 1925
$ws.Range("I119").Value = 1925
$ws.Range("K119").Value = 5775
$ws.Range("M119").Value = -937

$ws.Range("H132").Value = 1860.9375
$ws.Range("I132").Value = 2049.3333
$ws.Range("J132").Value = 1747.9
$ws.Range("K132").Value = 18443.9997
$ws.Range("L132").Value = 15731.1
$ws.Range("M132").Value = -15913.9997
$ws.Range("N132").Value = -20791.1

$ws.Range("H137").Value = 2057.5
$ws.Range("J137").Value = 2398.5
$ws.Range("L137").Value = 7195.5
$ws.Range("N137").Value = -17395.5

$ws.Range("H138").Value = 3100.0908
$ws.Range("I138").Value = 1665.4
$ws.Range("J138").Value = 4295.6665
$ws.Range("K138").Value = 4996.200000000001
$ws.Range("L138").Value = 12886.9995
$ws.Range("M138").Value = 143.7999999999993
$ws.Range("N138").Value = -23166.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H54").Value = 9167.833000000001
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 9167.833000000001
$ws.Range("K54").Value = 0
$ws.Range("L54").ClearContents()
$ws.Range("M54").Value = 9167.833000000001
$ws.Range("N54").Value = -9947.833000000001

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").ClearContents()
$ws.Range("N128").Value = 0

$ws.Range("H132").Value = 5090.421
$ws.Range("I132").Value = 2291.6333
$ws.Range("K132").Value = 6874.8999
$ws.Range("M132").Value = -4344.8999

$ws.Range("H141").Value = 144995
$ws.Range("J141").Value = 144995
$ws.Range("L141").Value = 144995
$ws.Range("N141").Value = -155355

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 6575.3125
$ws.Range("I93").Value = 5073.357
$ws.Range("J93").Value = 7743.5
$ws.Range("K93").Value = 5073.357
$ws.Range("L93").Value = 7743.5
$ws.Range("M93").Value = -3825.357
$ws.Range("N93").Value = -10239.5

$ws.Range("H132").Value = 1259508.1
$ws.Range("I132").Value = 2450.6875
$ws.Range("J132").Value = 2516565.5
$ws.Range("K132").Value = 7352.0625
$ws.Range("L132").Value = 7549696.5
$ws.Range("M132").Value = -4822.0625
$ws.Range("N132").Value = -7554756.5

$ws.Range("H136").Value = 18207
$ws.Range("I136").Value = 18202.309
$ws.Range("K136").Value = 54606.927
$ws.Range("M136").Value = -52056.927

$ws.Range("H140").Value = 178246.88
$ws.Range("J140").Value = 178246.88
$ws.Range("L140").Value = 178246.88
$ws.Range("N140").Value = -188606.88

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 73952.19
$ws.Range("J64").Value = 73952.19
$ws.Range("L64").Value = 73952.19
$ws.Range("N64").Value = -74448.19

$ws.Range("H67").Value = 73952.19
$ws.Range("J67").Value = 73952.19
$ws.Range("L67").Value = 73952.19
$ws.Range("N67").Value = -75668.19

$ws.Range("H140").Value = 138005.2
$ws.Range("J140").Value = 138005.2
$ws.Range("L140").Value = 138005.2
$ws.Range("N140").Value = -148365.2

$ws.Range("H141").Value = 70280
$ws.Range("J141").Value = 70280
$ws.Range("L141").Value = 70280
$ws.Range("N141").Value = -80640
